$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.116.76'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.917.52'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.62%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.22'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.34%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5064'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.57%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08353'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.21%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.84%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.03'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.79%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.914.24'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.49%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.418'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.75%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.246'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.41%  '

$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.002'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.55'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.03%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001097'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.64%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06514'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.13%  '

$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.53'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.81%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.03%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.945'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.71%  '

$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '30.127.08'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.69%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.194'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.06%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.128.07'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.16%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '21.86'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.62%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.95'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.272'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.31%  '

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '128.93'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.61%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.143'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.27%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1045'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.17%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.950'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.786'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.38%  '

$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.02456'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.91%  '

$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.362'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.84%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06447'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.90%  '

$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2157'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.02%  '

$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6524'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.82%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.198'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.27%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.613'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.53%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.40'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.63%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.211'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.70%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.37'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.66%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6077'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.40%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.193'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +10.69%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.624'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.37'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.210'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.71%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '79.03'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.49%  '

$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.137'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.30%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06821'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.36%  '
